# Apply the prediction_tracker.xlsx edits described by the commit diff.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Sheet "Predictions": drop the Expected/Base-Expected-Goals + Target
#    Score Input columns (Y:AC), update a few Actual_Score/Actual_Result
#    cells, revert the AS Monaco vs PSG "actual updated" row, replace the
#    former "Cagliari vs Lecce" 2026-02-18 row with the new
#    "Borussia Dortmund vs Atalanta" 2026-02-17 fixture, and drop the
#    trailing AS Monaco duplicate row (the sheet shrinks from 30 to 29
#    rows).
# ---------------------------------------------------------------------
$wsPred = $wb.Worksheets.Item("Predictions")

# Remove columns Y:AC (Expected_Goals_Home ... Target_Score_Input)
$wsPred.Range("Y1:AC30").Delete()

# Row 24: Lyon vs Nice now has a confirmed actual score
$wsPred.Range("K24").Value = "2-0"
$wsPred.Range("L24").Value = "Home"

# Row 26: Girona vs Barcelona now has a confirmed actual score
$wsPred.Range("K26").Value = "2-1"
$wsPred.Range("L26").Value = "Home"

# Row 27: Cagliari vs Lecce (2026-02-16) now has a confirmed actual score
$wsPred.Range("K27").Value = "0-2"
$wsPred.Range("L27").Value = "Away"

# Row 28: AS Monaco vs Paris Saint-Germain - revert the user-entered actual
# result back to blank/unverified.
$wsPred.Range("K28").Value = ""
$wsPred.Range("L28").Value = ""
$wsPred.Range("M28").Value = ""
$wsPred.Range("N28").Value = ""
$wsPred.Range("O28").Value = ""
$wsPred.Range("P28").Value = ""

# Row 29: replace "Cagliari vs Lecce" (2026-02-18) with the new fixture
# "Borussia Dortmund vs Atalanta" (2026-02-17, Champions_League).
# Force text formatting first so the date-like string isn't auto-converted
# into a real Excel date value.
$wsPred.Range("A29").NumberFormat = "@"
$wsPred.Range("A29").Value = "2026-02-17"
$wsPred.Range("B29").Value = "Borussia Dortmund vs Atalanta"
$wsPred.Range("C29").Value = "Champions_League"
$wsPred.Range("D29").Value = "Borussia Dortmund"
$wsPred.Range("E29").Value = "Atalanta"
$wsPred.Range("F29").Value = 37.87482918330498
$wsPred.Range("G29").Value = 29.19433780269773
$wsPred.Range("H29").Value = 32.93083301399729
$wsPred.Range("Q29").Value = "Borussia Dortmund vs Atalanta"
$wsPred.Range("R29").Value = "Borussia Dortmund"
$wsPred.Range("S29").Value = "Atalanta"
$wsPred.Range("W29").Value = "League mismatch from datasets: home=Bundesliga, away=Serie_A. Using 'Champions_League'."
$wsPred.Range("X29").Value = "match=Borussia Dortmund vs Atalanta; date=2026-02-17; league=Champions_League"

# Row 30 (old AS Monaco vs Paris Saint-Germain duplicate) is dropped entirely.
$wsPred.Rows.Item(30).Delete()

# ---------------------------------------------------------------------
# 2) Sheet "Summary": refreshed aggregate stats now that there are more
#    verified matches.
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B2").Value = 21
$wsSummary.Range("B3").Value = 11
$wsSummary.Range("B4").Value = 52.38
$wsSummary.Range("B6").Value = 19.05
$wsSummary.Range("B7").Value = 1.5

# ---------------------------------------------------------------------
# 3) Sheet "bet data": replace the "Cagliari vs Lecce" (2026-02-18) row
#    with the new Borussia Dortmund vs Atalanta fixture, and drop the
#    trailing AS Monaco duplicate row (the sheet shrinks from 12 to 11
#    rows).
# ---------------------------------------------------------------------
$wsBetData = $wb.Worksheets.Item("bet data")
$wsBetData.Range("A11").NumberFormat = "@"
$wsBetData.Range("A11").Value = "2026-02-17"
$wsBetData.Range("B11").Value = "Borussia Dortmund vs Atalanta"
$wsBetData.Range("C11").Value = "Champions_League"
$wsBetData.Range("D11").Value = 1.8
$wsBetData.Range("E11").Value = 6.3
$wsBetData.Range("F11").Value = 17.7
$wsBetData.Range("G11").Value = 39.1
$wsBetData.Range("H11").Value = 65.90000000000001
$wsBetData.Range("I11").Value = 85.59999999999999
$wsBetData.Range("J11").Value = 95.3
$wsBetData.Range("K11").Value = 92
$wsBetData.Range("L11").Value = 8
$wsBetData.Range("M11").Value = 71.90000000000001
$wsBetData.Range("N11").Value = 28.1
$wsBetData.Range("O11").Value = 46.4
$wsBetData.Range("P11").Value = 53.6
$wsBetData.Range("Q11").Value = 24.9
$wsBetData.Range("R11").Value = 75.09999999999999

$wsBetData.Rows.Item(12).Delete()

# ---------------------------------------------------------------------
# 4) Sheet "bet predic": a few bets now have recorded actual
#    scores/results, one bet's recorded result was reverted, the
#    "Cagliari vs Lecce" row becomes "Borussia Dortmund vs Atalanta",
#    and the trailing AS Monaco duplicate row is dropped (12 -> 11 rows).
# ---------------------------------------------------------------------
$wsBetPredic = $wb.Worksheets.Item("bet predic")

$wsBetPredic.Range("F6").Value = "2-0"
$wsBetPredic.Range("G6").Value = "Won"

$wsBetPredic.Range("F8").Value = "2-1"
$wsBetPredic.Range("G8").Value = "Won"

$wsBetPredic.Range("F9").Value = "0-2"
$wsBetPredic.Range("G9").Value = "Lost"

$wsBetPredic.Range("F10").Value = ""
$wsBetPredic.Range("G10").Value = ""

$wsBetPredic.Range("A11").NumberFormat = "@"
$wsBetPredic.Range("A11").Value = "2026-02-17"
$wsBetPredic.Range("B11").Value = "Borussia Dortmund vs Atalanta"
$wsBetPredic.Range("E11").Value = "Model-only selection by highest probability 0.820."
$wsBetPredic.Range("J11").Value = 0.8204

$wsBetPredic.Rows.Item(12).Delete()

# ---------------------------------------------------------------------
# 5) Drop the "Model Eval", "Model Eval League" and "Model Eval Segments"
#    sheets entirely.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("Model Eval Segments").Delete()
$wb.Worksheets.Item("Model Eval League").Delete()
$wb.Worksheets.Item("Model Eval").Delete()

# Restore the originally active sheet/tab.
$wsPred.Activate()
